# Update the six LinkedIn carousel slides: replace the NTPC Green Energy /
# GAIL joint-venture narrative with the Bhuj solar project expansion
# narrative, paragraph by paragraph, leaving paragraph formatting
# (a:pPr / a:defRPr) and the text box's size untouched.

function Set-ParagraphTexts($shape, $oldTexts, $newTexts) {
    $tr = $shape.TextFrame.TextRange

    # spAutoFit shapes can resize on text edits; remember the original
    # size so it can be restored once the new text is in place.
    $origHeight = $shape.Height
    $origWidth  = $shape.Width

    $count = $oldTexts.Count
    $starts = @()
    $lens   = @()

    $pos = 1
    for ($i = 0; $i -lt $count; $i++) {
        $len = $oldTexts[$i].Length
        $starts += $pos
        $lens += $len
        # Each paragraph break consumes one character position.
        $pos = $pos + $len + 1
    }

    # Replace from the last paragraph back to the first so that the
    # not-yet-processed offsets stay valid even though the lengths of
    # the new strings differ from the originals.
    for ($i = $count - 1; $i -ge 0; $i--) {
        $tr.Characters($starts[$i], $lens[$i]).Text = $newTexts[$i]
    }

    $shape.Height = $origHeight
    $shape.Width  = $origWidth
}

$p = $ppt.ActivePresentation

Set-ParagraphTexts $p.Slides.Item(1).Shapes.Item(1) `
    @("Joint Venture Formation", `
      "NTPC Green Energy and GAIL have formed a joint venture.", `
      "The joint venture operates on a 50:50 ownership basis.") `
    @("Project Expansion Overview", `
      "NTPC Green Energy has expanded the Bhuj solar project.", `
      "An additional capacity of 37.5 MW has been added.")

Set-ParagraphTexts $p.Slides.Item(2).Shapes.Item(1) `
    @("Focus on Renewable Energy", `
      "The joint venture is specifically focused on renewable energy.", `
      "This collaboration aims to enhance sustainable energy production.") `
    @("Total Capacity Achieved", `
      "The total capacity of the Bhuj solar project now reaches 8,347.78 MW.", `
      "This expansion contributes to NTPC Green Energy's overall capacity growth.")

Set-ParagraphTexts $p.Slides.Item(3).Shapes.Item(1) `
    @("Participants Overview", `
      "NTPC Green Energy is a subsidiary of NTPC Limited.", `
      "GAIL is a major player in the natural gas sector in India.") `
    @("Location of the Project", `
      "The expanded solar project is located in Bhuj.", `
      "Bhuj is a strategic location for solar energy production.")

Set-ParagraphTexts $p.Slides.Item(4).Shapes.Item(1) `
    @("Strategic Goals", `
      "The joint venture aligns with India's goals for renewable energy expansion.", `
      "It supports the transition towards cleaner energy sources.") `
    @("Company Profile", `
      "NTPC Green Energy focuses on renewable energy projects.", `
      "The company is part of NTPC Limited, a major player in India's energy sector.")

Set-ParagraphTexts $p.Slides.Item(5).Shapes.Item(1) `
    @("Project Development", `
      "Details on specific projects from the joint venture were not provided.", `
      "Future announcements regarding project specifics are anticipated.") `
    @("Importance of Solar Energy", `
      "The addition of renewable energy sources is critical for sustainability.", `
      "Solar energy helps reduce carbon emissions and dependence on fossil fuels.")

Set-ParagraphTexts $p.Slides.Item(6).Shapes.Item(1) `
    @("Industry Impact", `
      "The collaboration is expected to contribute to the growth of the renewable energy sector in India.", `
      "It signifies a partnership between two significant entities in the energy market.") `
    @("Future Outlook", `
      "NTPC Green Energy aims to further increase its renewable energy capacity.", `
      "Continued investment in solar projects is expected to support growth.")
